$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1285113333333333
$ws.Range("H2").Value = 0.385534
$ws.Range("I2").Value = 0.03749201237720504
$ws.Range("J2").Value = 0.03749201237720504
$ws.Range("M2").Value = 247.0944516666667
$ws.Range("N2").Value = 741.283355
$ws.Range("O2").Value = 0.8050739182622993
$ws.Range("P2").Value = 0.8050739182622993
$ws.Range("Q2").Value = 31.75443744295223
$ws.Range("R2").Value = 285.78993698657
$ws.Range("S2").Value = 0.03018384130805508
$ws.Range("T2").Value = 0.03018384130805509
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1285113333333333
$ws.Range("H3").Value = 0.385534
$ws.Range("I3").Value = 0.03749201237720504
$ws.Range("J3").Value = 0.03749201237720504
$ws.Range("O3").Value = 0.1379009747488701
$ws.Range("P3").Value = 0.13790097474887
$ws.Range("Q3").Value = 5.439212197355556
$ws.Range("R3").Value = 48.95290977619999
$ws.Range("S3").Value = 0.005170185052113275
$ws.Range("T3").Value = 0.005170185052113275
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1285113333333333
$ws.Range("H4").Value = 0.385534
$ws.Range("I4").Value = 0.03749201237720504
$ws.Range("J4").Value = 0.03749201237720504
$ws.Range("M4").Value = 11.590146
$ws.Range("N4").Value = 34.770438
$ws.Range("O4").Value = 0.03776258103132013
$ws.Range("P4").Value = 0.03776258103132013
$ws.Range("Q4").Value = 1.489465115988
$ws.Range("R4").Value = 13.405186043892
$ws.Range("S4").Value = 0.001415795155421462
$ws.Range("T4").Value = 0.001415795155421463
# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1285113333333333
$ws.Range("H5").Value = 0.385534
$ws.Range("I5").Value = 0.03749201237720504
$ws.Range("J5").Value = 0.03749201237720504
$ws.Range("M5").Value = 5.912082333333333
$ws.Range("N5").Value = 17.736247
$ws.Range("O5").Value = 0.01926252595751047
$ws.Range("P5").Value = 0.01926252595751047
$ws.Range("Q5").Value = 0.7597695834331111
$ws.Range("R5").Value = 6.837926250898
$ws.Range("S5").Value = 0.000722190861615216
$ws.Range("T5").Value = 0.0007221908616152161
# Row 6
$ws.Range("I6").Value = 0.7552862722193517
$ws.Range("J6").Value = 0.755286272219352
$ws.Range("M6").Value = 247.0944516666667
$ws.Range("N6").Value = 741.283355
$ws.Range("O6").Value = 0.8050739182622993
$ws.Range("P6").Value = 0.8050739182622993
$ws.Range("Q6").Value = 639.7013433531233
$ws.Range("R6").Value = 5757.312090178109
$ws.Range("S6").Value = 0.6080612785853592
$ws.Range("T6").Value = 0.6080612785853593
# Row 7
$ws.Range("I7").Value = 0.7552862722193517
$ws.Range("J7").Value = 0.755286272219352
$ws.Range("O7").Value = 0.1379009747488701
$ws.Range("P7").Value = 0.13790097474887
$ws.Range("S7").Value = 0.104154713153489
$ws.Range("T7").Value = 0.104154713153489
# Row 8
$ws.Range("I8").Value = 0.7552862722193517
$ws.Range("J8").Value = 0.755286272219352
$ws.Range("M8").Value = 11.590146
$ws.Range("N8").Value = 34.770438
$ws.Range("O8").Value = 0.03776258103132013
$ws.Range("P8").Value = 0.03776258103132013
$ws.Range("Q8").Value = 30.005659438524
$ws.Range("R8").Value = 270.050934946716
$ws.Range("S8").Value = 0.02852155905652698
$ws.Range("T8").Value = 0.02852155905652699
# Row 9
$ws.Range("I9").Value = 0.7552862722193517
$ws.Range("J9").Value = 0.755286272219352
$ws.Range("M9").Value = 5.912082333333333
$ws.Range("N9").Value = 17.736247
$ws.Range("O9").Value = 0.01926252595751047
$ws.Range("P9").Value = 0.01926252595751047
$ws.Range("Q9").Value = 15.30575448027266
$ws.Range("R9").Value = 137.751790322454
$ws.Range("S9").Value = 0.01454872142397658
$ws.Range("T9").Value = 0.01454872142397659
# Row 10
$ws.Range("G10").Value = 0.692415
$ws.Range("H10").Value = 2.077245
$ws.Range("I10").Value = 0.2020057770533527
$ws.Range("J10").Value = 0.2020057770533527
$ws.Range("M10").Value = 247.0944516666667
$ws.Range("N10").Value = 741.283355
$ws.Range("O10").Value = 0.8050739182622993
$ws.Range("P10").Value = 0.8050739182622993
$ws.Range("Q10").Value = 171.091904750775
$ws.Range("R10").Value = 1539.827142756975
$ws.Range("S10").Value = 0.1626295824439631
$ws.Range("T10").Value = 0.1626295824439632
# Row 11
$ws.Range("G11").Value = 0.692415
$ws.Range("H11").Value = 2.077245
$ws.Range("I11").Value = 0.2020057770533527
$ws.Range("J11").Value = 0.2020057770533527
$ws.Range("O11").Value = 0.1379009747488701
$ws.Range("P11").Value = 0.13790097474887
$ws.Range("Q11").Value = 29.3063033115
$ws.Range("R11").Value = 263.7567298035
$ws.Range("S11").Value = 0.02785679356056027
$ws.Range("T11").Value = 0.02785679356056027
# Row 12
$ws.Range("G12").Value = 0.692415
$ws.Range("H12").Value = 2.077245
$ws.Range("I12").Value = 0.2020057770533527
$ws.Range("J12").Value = 0.2020057770533527
$ws.Range("M12").Value = 11.590146
$ws.Range("N12").Value = 34.770438
$ws.Range("O12").Value = 0.03776258103132013
$ws.Range("P12").Value = 0.03776258103132013
$ws.Range("Q12").Value = 8.02519094259
$ws.Range("R12").Value = 72.22671848331
$ws.Range("S12").Value = 0.00762825952477202
$ws.Range("T12").Value = 0.007628259524772021
# Row 13
$ws.Range("G13").Value = 0.692415
$ws.Range("H13").Value = 2.077245
$ws.Range("I13").Value = 0.2020057770533527
$ws.Range("J13").Value = 0.2020057770533527
$ws.Range("M13").Value = 5.912082333333333
$ws.Range("N13").Value = 17.736247
$ws.Range("O13").Value = 0.01926252595751047
$ws.Range("P13").Value = 0.01926252595751047
$ws.Range("Q13").Value = 4.093614488835
$ws.Range("R13").Value = 36.842530399515
$ws.Range("S13").Value = 0.00389114152405728
$ws.Range("T13").Value = 0.003891141524057281
# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01787866666666667
$ws.Range("H14").Value = 0.053636
$ws.Range("I14").Value = 0.005215938350090445
$ws.Range("J14").Value = 0.005215938350090446
$ws.Range("M14").Value = 247.0944516666667
$ws.Range("N14").Value = 741.283355
$ws.Range("O14").Value = 0.8050739182622993
$ws.Range("P14").Value = 0.8050739182622993
$ws.Range("Q14").Value = 4.417719336531111
$ws.Range("R14").Value = 39.75947402878001
$ws.Range("S14").Value = 0.004199215924921907
$ws.Range("T14").Value = 0.004199215924921908
# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01787866666666667
$ws.Range("H15").Value = 0.053636
$ws.Range("I15").Value = 0.005215938350090445
$ws.Range("J15").Value = 0.005215938350090446
$ws.Range("O15").Value = 0.1379009747488701
$ws.Range("P15").Value = 0.13790097474887
$ws.Range("Q15").Value = 0.7567103949777778
$ws.Range("R15").Value = 6.8103935548
$ws.Range("S15").Value = 0.0007192829827074854
$ws.Range("T15").Value = 0.0007192829827074854
# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01787866666666667
$ws.Range("H16").Value = 0.053636
$ws.Range("I16").Value = 0.005215938350090445
$ws.Range("J16").Value = 0.005215938350090446
$ws.Range("M16").Value = 11.590146
$ws.Range("N16").Value = 34.770438
$ws.Range("O16").Value = 0.03776258103132013
$ws.Range("P16").Value = 0.03776258103132013
$ws.Range("Q16").Value = 0.207216356952
$ws.Range("R16").Value = 1.864947212568
$ws.Range("S16").Value = 0.0001969672945996607
$ws.Range("T16").Value = 0.0001969672945996607
# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01787866666666667
$ws.Range("H17").Value = 0.053636
$ws.Range("I17").Value = 0.005215938350090445
$ws.Range("J17").Value = 0.005215938350090446
$ws.Range("M17").Value = 5.912082333333333
$ws.Range("N17").Value = 17.736247
$ws.Range("O17").Value = 0.01926252595751047
$ws.Range("Q17").Value = 0.1057001493435556
$ws.Range("R17").Value = 0.951301344092
$ws.Range("S17").Value = 0.0001004721478613915
$ws.Range("T17").Value = 0.0001004721478613916
